# Add ECS wave support: insert a new row for "ECS last wave" right after
# the existing "ECS sheet" row on the Parameters sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Parameters")

# Insert a new row at row 3 (pushes VPC/Subnet/Secgroup/EnterpriseProject rows down by one)
$ws.Rows.Item(3).Insert([Microsoft.Office.Interop.Excel.XlInsertShiftDirection]::xlShiftDown, [Microsoft.Office.Interop.Excel.XlInsertFormatOrigin]::xlFormatFromLeftOrAbove)

# Copy the row-above formatting into the newly inserted row (matches Excel's
# default "Insert Copied Cells"/format-from-above behaviour)
$ws.Range("A2:B2").Copy()
$ws.Range("A3:B3").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

# Populate the newly inserted row
$ws.Cells.Item(3, 1).Value = "ECS last wave"
$ws.Cells.Item(3, 2).Value = 1

# Reset the selection back to the top-left cell (clears the stale B8 selection)
$ws.Range("A1").Select()
